$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: column A used to read "National" / column B "Unemployment rate"
# Now column A header becomes "State"
$ws.Range("A1").Value = "State"

# Row 2 column A used to read "United States" (the national total row),
# it is now labeled "National"
$ws.Range("A2").Value = "National"

# Update the active cell/selection to E2 to match the saved view state
$ws.Range("E2").Select()
